# Update cfb_weather.xlsx with Timestamp 2024-10-22T16:21:45.364535
#
# This script:
#  1. Updates the run Timestamp (column AK on "FBS") for every game row to
#     the new scrape time.
#  2. Applies the scraped odds/wind-direction corrections that came in with
#     that run, on both the "FBS" and "Other" sheets.

$wb = $excel.ActiveWorkbook

$fbs   = $wb.Worksheets.Item("FBS")
$other = $wb.Worksheets.Item("Other")

$newTimestamp = "2024-10-22T16:21:45.364535"

# ---------------------------------------------------------------------
# 1) Timestamp refresh - column AK, rows 2 through 54 on "FBS"
# ---------------------------------------------------------------------
for ($r = 2; $r -le 54; $r++) {
    $fbs.Cells.Item($r, 37).Value = $newTimestamp
}

# ---------------------------------------------------------------------
# 2) FBS sheet corrections
# ---------------------------------------------------------------------

# Row 4 - Central Michigan @ Miami (OH): Open/Current line pulled, clear
$fbs.Range("AA4").Value = $null
$fbs.Range("AB4").Value = $null
$fbs.Range("AF4").Value = $null

# Row 6 - Texas Tech @ TCU
$fbs.Range("AB6").Value = -6.5
$fbs.Range("AF6").Value = 0

# Row 9 - Arkansas @ Mississippi State
$fbs.Range("AB9").Value = 6.5
$fbs.Range("AF9").Value = 0.5

# Row 13 - Sam Houston State @ Florida International
$fbs.Range("Q13").Value = "S"
$fbs.Range("Z13").Value = -108

# Row 15 - Liberty @ Kennesaw State
$fbs.Range("AB15").Value = 25.5
$fbs.Range("AF15").Value = -2

# Row 19 - Louisville @ Boston College
$fbs.Range("Q19").Value = "NNW"

# Row 20 - Rutgers @ USC
$fbs.Range("AB20").Value = -14
$fbs.Range("AF20").Value = -0.5

# Row 30 - Brigham Young @ UCF
$fbs.Range("AB30").Value = -1
$fbs.Range("AF30").Value = 2.5

# Row 32 - Southern Miss @ James Madison
$fbs.Range("AB32").Value = -24
$fbs.Range("AF32").Value = 0

# Row 37 - Wake Forest @ Stanford
$fbs.Range("Q37").Value = "S"

# Row 39 - Illinois @ Oregon
$fbs.Range("Q39").Value = "E"

# Row 45 - Michigan State @ Michigan
$fbs.Range("AB45").Value = -4
$fbs.Range("AF45").Value = -2

# Row 50 - Kansas @ Kansas State
$fbs.Range("AB50").Value = -10.5
$fbs.Range("AF50").Value = 0.5

# Row 53 - Washington State @ San Diego State
$fbs.Range("Z53").Value = -106

# ---------------------------------------------------------------------
# 3) Other sheet corrections (wind_dir_fg, column S)
# ---------------------------------------------------------------------

# Row 14 - East Tennessee State vs Wofford
$other.Range("S14").Value = "E"

# Row 22 - Cornell vs Brown
$other.Range("S22").Value = "E"

# Row 23 - Richmond vs Bryant University
$other.Range("S23").Value = "E"
